# Update "想去人数" (number of people wanting to attend) figures on the
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets
# to reflect the regenerated site data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 872
$ws.Range("F10").Value = 34
$ws.Range("F11").Value = 2695
$ws.Range("F13").Value = 1572
$ws.Range("F14").Value = 7246
$ws.Range("F16").Value = 7390
$ws.Range("F17").Value = 9
$ws.Range("F18").Value = 24
$ws.Range("F19").Value = 5768
$ws.Range("F20").Value = 3161
$ws.Range("F21").Value = 3536
$ws.Range("F25").Value = 209
$ws.Range("F26").Value = 1987
$ws.Range("F28").Value = 323
$ws.Range("F29").Value = 898
$ws.Range("F30").Value = 241
$ws.Range("F31").Value = 713
$ws.Range("F33").Value = 2507
$ws.Range("F34").Value = 1322
$ws.Range("F35").Value = 2969
$ws.Range("F36").Value = 96
$ws.Range("F38").Value = 184
$ws.Range("F39").Value = 434
$ws.Range("F40").Value = 1155

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 381

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 872
$ws.Range("F12").Value = 34
$ws.Range("F14").Value = 2695
$ws.Range("F15").Value = 1572
$ws.Range("F17").Value = 7246
$ws.Range("F19").Value = 7390
$ws.Range("F20").Value = 9
$ws.Range("F21").Value = 24
$ws.Range("F22").Value = 5768
$ws.Range("F23").Value = 3161
$ws.Range("F24").Value = 3536
$ws.Range("F30").Value = 1987
$ws.Range("F33").Value = 323
$ws.Range("F34").Value = 898
$ws.Range("F35").Value = 241
$ws.Range("F36").Value = 713
$ws.Range("F38").Value = 2507
$ws.Range("F39").Value = 1322
$ws.Range("F41").Value = 2969
$ws.Range("F42").Value = 96
$ws.Range("F44").Value = 184
$ws.Range("F46").Value = 434
$ws.Range("F47").Value = 1155
